$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scene")
if (-not $ws) { $ws = $wb.ActiveSheet }

$ws.Range("D6").Value = 19
$ws.Range("D8").Value = 5
$ws.Range("D9").Value = 8
$ws.Range("D11").Value = 7
$ws.Range("D12").Value = 15
$ws.Range("D13").Value = 4
$ws.Range("D14").Value = 10
$ws.Range("D15").Value = 12
$ws.Range("D16").Value = 3
$ws.Range("D17").Value = 6
$ws.Range("D18").Value = 14
$ws.Range("D19").Value = 16
$ws.Range("D20").Value = 16
$ws.Range("D21").Value = 18
$ws.Range("D22").Value = 20

$ws.Range("D6").Select()
